$wb = $excel.ActiveWorkbook

# Rename the regression sheets to the "PL_" prefixed names used in the
# working version of the PL model.
$wsR1a = $wb.Worksheets.Item("R1a")
$wsR1a.Name = "PL_R1a"

$wsR1b = $wb.Worksheets.Item("R1b")
$wsR1b.Name = "PL_R1b"

# The active tab moves from PL_R1a to PL_R1b.
$wsR1b.Activate()
